# New crime data collected - weekly CompStat update for the 45th Precinct.
# Moves the report window forward one week (Volume/Number + date range)
# and refreshes the crime-complaint figures for the new reporting week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: issue number and the two report-covering dates.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/6/2023  Through  11/12/2023"

# ---------------------------------------------------------------------
# A few cells flip between "numeric count" and the report's textual
# placeholders ("0" / "***.*"). Copying a same-shaped neighbor cell
# keeps the original number format / style pooling intact.
# ---------------------------------------------------------------------

# Rape/UCR Rape*/Other Sex Crimes week-to-date 2023 count becomes the
# placeholder text "0" -- copy the whole cell (value + format) from the
# neighboring "2022" column, which already holds that exact placeholder.
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("D26").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Range("D27").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Transit / Housing / Hate Crimes rows: some placeholder text cells turn
# into real numeric counts this week. Copy number-format only from a
# same-row cell already in the right numeric style, then write the value.
$ws.Range("F22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("F22").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("K22").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("K22").Copy()
$ws.Range("H22").PasteSpecial(-4122)

$ws.Range("F23").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("F23").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("H23").Copy()
$ws.Range("E23").PasteSpecial(-4122)

$ws.Range("G30").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("H30").Copy()
$ws.Range("E30").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = -80
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 26
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 177
$ws.Range("J16").Value = 187
$ws.Range("K16").Value = -5.347593582887
$ws.Range("L16").Value = 17.218543046357
$ws.Range("M16").Value = -12.807881773399
$ws.Range("N16").Value = -55.970149253731

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -77.777777777777
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -50
$ws.Range("I17").Value = 254
$ws.Range("J17").Value = 231
$ws.Range("K17").Value = 9.956709956709
$ws.Range("L17").Value = 25.123152709359
$ws.Range("M17").Value = 57.763975155279
$ws.Range("N17").Value = 26.368159203980

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 25
$ws.Range("I18").Value = 121
$ws.Range("J18").Value = 90
$ws.Range("K18").Value = 34.444444444444
$ws.Range("L18").Value = 24.742268041237
$ws.Range("M18").Value = -55.185185185185
$ws.Range("N18").Value = -83.356258596973

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 66
$ws.Range("G19").Value = 54
$ws.Range("H19").Value = 22.222222222222
$ws.Range("I19").Value = 580
$ws.Range("J19").Value = 498
$ws.Range("K19").Value = 16.465863453815
$ws.Range("L19").Value = 41.463414634146
$ws.Range("M19").Value = 38.424821002386
$ws.Range("N19").Value = 52.631578947368

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 11
$ws.Range("E20").Value = 83.333333333333
$ws.Range("F20").Value = 37
$ws.Range("H20").Value = 27.586206896551
$ws.Range("I20").Value = 456
$ws.Range("J20").Value = 264
$ws.Range("K20").Value = 72.727272727272
$ws.Range("L20").Value = 70.149253731343
$ws.Range("M20").Value = 150.549450549451
$ws.Range("N20").Value = -75.497044599677

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 43
$ws.Range("E21").Value = -23.255813953488
$ws.Range("F21").Value = 141
$ws.Range("G21").Value = 139
$ws.Range("H21").Value = 1.438848920863
$ws.Range("I21").Value = 1609
$ws.Range("J21").Value = 1286
$ws.Range("K21").Value = 25.116640746500
$ws.Range("L21").Value = 40.401396160558
$ws.Range("M21").Value = 28.411811652035
$ws.Range("N21").Value = -55.355160932297

# ---------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 400
$ws.Range("I22").Value = 17
$ws.Range("J22").Value = 14
$ws.Range("K22").Value = 21.428571428571
$ws.Range("L22").Value = 30.769230769230
$ws.Range("M22").Value = 41.666666666666

# ---------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 54
$ws.Range("J23").Value = 29
$ws.Range("K23").Value = 86.206896551724
$ws.Range("L23").Value = 28.571428571428
$ws.Range("M23").Value = 14.893617021276

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = 33.333333333333
$ws.Range("F24").Value = 86
$ws.Range("G24").Value = 90
$ws.Range("H24").Value = -4.444444444444
$ws.Range("I24").Value = 1210
$ws.Range("J24").Value = 1006
$ws.Range("K24").Value = 20.278330019880
$ws.Range("L24").Value = 41.686182669789
$ws.Range("M24").Value = -6.055900621118

# ---------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 42
$ws.Range("G25").Value = 43
$ws.Range("H25").Value = -2.325581395348
$ws.Range("I25").Value = 446
$ws.Range("J25").Value = 402
$ws.Range("K25").Value = 10.945273631840
$ws.Range("L25").Value = 17.368421052631
$ws.Range("M25").Value = 18.617021276595

# ---------------------------------------------------------------------
# Row 27 - Other Sex Crimes (remaining numeric cells)
# ---------------------------------------------------------------------
$ws.Range("F27").Value = 5
$ws.Range("H27").Value = 66.666666666666
$ws.Range("L27").Value = 2.631578947368

# ---------------------------------------------------------------------
# Row 28 - Shooting Vic.
# ---------------------------------------------------------------------
$ws.Range("L28").Value = -30.769230769230

# ---------------------------------------------------------------------
# Row 29 - Shooting Inc.
# ---------------------------------------------------------------------
$ws.Range("L29").Value = -30

# ---------------------------------------------------------------------
# Row 30 - Hate Crimes
# ---------------------------------------------------------------------
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = -100
$ws.Range("G30").Value = 2
$ws.Range("J30").Value = 6
$ws.Range("K30").Value = -33.333333333333
